$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Ccl3/Ackr2 -> FAPs) - values updated, labels unchanged
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.07287233333333333
$ws.Range("H2").Value = 0.218617
$ws.Range("I2").Value = 0.0006119733110023554
$ws.Range("J2").Value = 0.0006119733110023554
$ws.Range("Q2").Value = 0.03985305321355555
$ws.Range("R2").Value = 0.358677478922
$ws.Range("S2").Value = 0.0006119733110023554
$ws.Range("T2").Value = 0.0006119733110023554

# Row 3 now represents MuSCs -> Ccl3/Ackr2 -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.135447
$ws.Range("H3").Value = 0.406341
$ws.Range("I3").Value = 0.001137468024746511
$ws.Range("J3").Value = 0.001137468024746511
$ws.Range("Q3").Value = 0.07407442923400001
$ws.Range("R3").Value = 0.666669863106
$ws.Range("S3").Value = 0.001137468024746511
$ws.Range("T3").Value = 0.001137468024746511

# New row 4: Resolving-Mac -> Ccl3/Ackr2 -> FAPs
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Ccl3"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 118.8693136666667
$ws.Range("H4").Value = 356.607941
$ws.Range("I4").Value = 0.9982505586642512
$ws.Range("J4").Value = 0.9982505586642512
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5468886666666667
$ws.Range("N4").Value = 1.640666
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 65.00828045874511
$ws.Range("R4").Value = 585.074524128706
$ws.Range("S4").Value = 0.9982505586642512
$ws.Range("T4").Value = 0.9982505586642512
